$wb = $excel.ActiveWorkbook

# Step 1: Insert a new worksheet "num" as the first sheet (Worksheets.Add() with no args inserts at front).
$numSheet = $wb.Worksheets.Add()
$numSheet.Name = "num"

# Step 2: Translate / update the shared label text (same cells, new text) across every sheet.
# This rewrites B1/C1/D1 header cells and A2:A7 row-label cells on every OLD sheet (perc..rcr);
# for the brand new "num" sheet we will write these labels explicitly below together with its data.
$labelMap = @{
  "all" = "TUTTI GLI AUTORI";
  "italy" = "ITALIA";
  "fism" = "FISM";
  "any" = "tutte le pubblicazioni";
  "prevenzione_primaria" = "prevenzione primaria";
  "prevenzione_secondaria" = "prevenzione secondaria";
  "prevenzione_terziaria" = "prevenzione terziaria";
  "covid" = "covid";
  "altre_malattie" = "altre malattie";
}

$oldSheetNames = @("perc","if","altmetric","cima_index","cit","reccit","fcr","rcr")
foreach ($name in $oldSheetNames) {
    $ws = $wb.Worksheets.Item($name)
    foreach ($addr in @("B1","C1","D1")) {
        $cell = $ws.Range($addr)
        $t = $cell.Text
        if ($labelMap.ContainsKey($t)) {
            $cell.Value = $labelMap[$t]
        }
    }
    foreach ($r in 2..7) {
        $cell = $ws.Cells.Item($r, 1)
        $t = $cell.Text
        if ($labelMap.ContainsKey($t)) {
            $cell.Value = $labelMap[$t]
        }
    }
}

# Step 3: Write header labels (row 1) and row labels (column A) on the new "num" sheet, matching the other sheets.
$numSheet.Range("B1").Value = "TUTTI GLI AUTORI"
$numSheet.Range("C1").Value = "ITALIA"
$numSheet.Range("D1").Value = "FISM"
$numRowLabels = @("tutte le pubblicazioni","prevenzione primaria","prevenzione secondaria","prevenzione terziaria","covid","altre malattie")
for ($i = 0; $i -lt 6; $i++) {
    $numSheet.Cells.Item($i + 2, 1).Value = $numRowLabels[$i]
}

# Step 4: Write the final data values (B2:D7) for every sheet, in final left-to-right order.

# Sheet "num"
$ws = $wb.Worksheets.Item("num")
$ws.Cells.Item(2, 2).Value = 10767
$ws.Cells.Item(2, 3).Value = 1400
$ws.Cells.Item(2, 4).Value = 228
$ws.Cells.Item(3, 2).Value = 1021
$ws.Cells.Item(3, 3).Value = 132
$ws.Cells.Item(3, 4).Value = 25
$ws.Cells.Item(4, 2).Value = 973
$ws.Cells.Item(4, 3).Value = 172
$ws.Cells.Item(4, 4).Value = 41
$ws.Cells.Item(5, 2).Value = 1580
$ws.Cells.Item(5, 3).Value = 191
$ws.Cells.Item(5, 4).Value = 47
$ws.Cells.Item(6, 2).Value = 543
$ws.Cells.Item(6, 3).Value = 117
$ws.Cells.Item(6, 4).Value = 25
$ws.Cells.Item(7, 2).Value = 337
$ws.Cells.Item(7, 3).Value = 37
$ws.Cells.Item(7, 4).Value = 5

# Sheet "perc"
$ws = $wb.Worksheets.Item("perc")
$ws.Cells.Item(2, 2).Value = 100
$ws.Cells.Item(2, 3).Value = 100
$ws.Cells.Item(2, 4).Value = 100
$ws.Cells.Item(3, 2).Value = 9.482678554843503
$ws.Cells.Item(3, 3).Value = 9.428571428571429
$ws.Cells.Item(3, 4).Value = 10.96491228070175
$ws.Cells.Item(4, 2).Value = 9.036871923469862
$ws.Cells.Item(4, 3).Value = 12.28571428571429
$ws.Cells.Item(4, 4).Value = 17.98245614035088
$ws.Cells.Item(5, 2).Value = 14.67446828271571
$ws.Cells.Item(5, 3).Value = 13.64285714285714
$ws.Cells.Item(5, 4).Value = 20.6140350877193
$ws.Cells.Item(6, 2).Value = 5.043187517414322
$ws.Cells.Item(6, 3).Value = 8.357142857142858
$ws.Cells.Item(6, 4).Value = 10.96491228070175
$ws.Cells.Item(7, 2).Value = 3.129934057769109
$ws.Cells.Item(7, 3).Value = 2.642857142857143
$ws.Cells.Item(7, 4).Value = 2.192982456140351

# Sheet "if"
$ws = $wb.Worksheets.Item("if")
$ws.Cells.Item(2, 2).Value = 6.040066644113667
$ws.Cells.Item(2, 3).Value = 6.149315499606608
$ws.Cells.Item(2, 4).Value = 8.419077981651375
$ws.Cells.Item(3, 2).Value = 6.422116788321167
$ws.Cells.Item(3, 3).Value = 6.506359374999999
$ws.Cells.Item(3, 4).Value = 10.64583333333333
$ws.Cells.Item(4, 2).Value = 7.320967069154775
$ws.Cells.Item(4, 3).Value = 9.149571428571429
$ws.Cells.Item(4, 4).Value = 13.95365853658537
$ws.Cells.Item(5, 2).Value = 4.375853061224491
$ws.Cells.Item(5, 3).Value = 4.963585635359116
$ws.Cells.Item(5, 4).Value = 6.106744680851063
$ws.Cells.Item(6, 2).Value = 5.57431906614786
$ws.Cells.Item(6, 3).Value = 5.221290598290597
$ws.Cells.Item(6, 4).Value = 6.056
$ws.Cells.Item(7, 2).Value = 5.564712460063898
$ws.Cells.Item(7, 3).Value = 6.597142857142857
$ws.Cells.Item(7, 4).Value = 5.475

# Sheet "altmetric"
$ws = $wb.Worksheets.Item("altmetric")
$ws.Cells.Item(2, 2).Value = 14.06434115776676
$ws.Cells.Item(2, 3).Value = 13.45746962115797
$ws.Cells.Item(2, 4).Value = 31.02631578947368
$ws.Cells.Item(3, 2).Value = 33.59050445103858
$ws.Cells.Item(3, 3).Value = 15.71969696969697
$ws.Cells.Item(3, 4).Value = 20.64
$ws.Cells.Item(4, 2).Value = 19.76421923474664
$ws.Cells.Item(4, 3).Value = 38.48255813953488
$ws.Cells.Item(4, 4).Value = 106.4390243902439
$ws.Cells.Item(5, 2).Value = 7.948979591836735
$ws.Cells.Item(5, 3).Value = 8.705263157894738
$ws.Cells.Item(5, 4).Value = 10.06382978723404
$ws.Cells.Item(6, 2).Value = 24.42329020332717
$ws.Cells.Item(6, 3).Value = 11.35042735042735
$ws.Cells.Item(6, 4).Value = 16.04
$ws.Cells.Item(7, 2).Value = 21.55988023952096
$ws.Cells.Item(7, 3).Value = 13.61111111111111
$ws.Cells.Item(7, 4).Value = 5.4

# Sheet "cima_index"
$ws = $wb.Worksheets.Item("cima_index")
$ws.Cells.Item(2, 2).Value = 1.575926177228932
$ws.Cells.Item(2, 3).Value = 2.086519114688129
$ws.Cells.Item(2, 4).Value = 2.352941176470588
$ws.Cells.Item(3, 2).Value = 1.589211618257261
$ws.Cells.Item(3, 3).Value = 1.99009900990099
$ws.Cells.Item(3, 4).Value = 2.047619047619047
$ws.Cells.Item(4, 2).Value = 1.769911504424779
$ws.Cells.Item(4, 3).Value = 2.384
$ws.Cells.Item(4, 4).Value = 2.966666666666667
$ws.Cells.Item(5, 2).Value = 1.470691163604549
$ws.Cells.Item(5, 3).Value = 2.13768115942029
$ws.Cells.Item(5, 4).Value = 2.696969696969697
$ws.Cells.Item(6, 2).Value = 1.638613861386139
$ws.Cells.Item(6, 3).Value = 2.293478260869565
$ws.Cells.Item(6, 4).Value = 3.05
$ws.Cells.Item(7, 2).Value = 1.521951219512195
$ws.Cells.Item(7, 3).Value = 2
$ws.Cells.Item(7, 4).Value = 1

# Sheet "cit"
$ws = $wb.Worksheets.Item("cit")
$ws.Cells.Item(2, 2).Value = 7.839147105583092
$ws.Cells.Item(2, 3).Value = 10.56540385989993
$ws.Cells.Item(2, 4).Value = 14.47368421052632
$ws.Cells.Item(3, 2).Value = 10.94658753709199
$ws.Cells.Item(3, 3).Value = 14.79545454545454
$ws.Cells.Item(3, 4).Value = 21.08
$ws.Cells.Item(4, 2).Value = 11.19441571871768
$ws.Cells.Item(4, 3).Value = 17.47093023255814
$ws.Cells.Item(4, 4).Value = 12.78048780487805
$ws.Cells.Item(5, 2).Value = 7.260204081632653
$ws.Cells.Item(5, 3).Value = 9.726315789473684
$ws.Cells.Item(5, 4).Value = 11.78723404255319
$ws.Cells.Item(6, 2).Value = 16.75415896487985
$ws.Cells.Item(6, 3).Value = 21.90598290598291
$ws.Cells.Item(6, 4).Value = 41.68
$ws.Cells.Item(7, 2).Value = 9.482035928143713
$ws.Cells.Item(7, 3).Value = 13.86111111111111
$ws.Cells.Item(7, 4).Value = 11

# Sheet "reccit"
$ws = $wb.Worksheets.Item("reccit")
$ws.Cells.Item(2, 2).Value = 7.133358271766577
$ws.Cells.Item(2, 3).Value = 9.441029306647605
$ws.Cells.Item(2, 4).Value = 12.67105263157895
$ws.Cells.Item(3, 2).Value = 10.07715133531157
$ws.Cells.Item(3, 3).Value = 13.12878787878788
$ws.Cells.Item(3, 4).Value = 19.36
$ws.Cells.Item(4, 2).Value = 10.22854188210962
$ws.Cells.Item(4, 3).Value = 16.04651162790698
$ws.Cells.Item(4, 4).Value = 12
$ws.Cells.Item(5, 2).Value = 6.48405612244898
$ws.Cells.Item(5, 3).Value = 8.58421052631579
$ws.Cells.Item(5, 4).Value = 10.06382978723404
$ws.Cells.Item(6, 2).Value = 13.95933456561922
$ws.Cells.Item(6, 3).Value = 17.68376068376068
$ws.Cells.Item(6, 4).Value = 31.24
$ws.Cells.Item(7, 2).Value = 8.865269461077844
$ws.Cells.Item(7, 3).Value = 12.91666666666667
$ws.Cells.Item(7, 4).Value = 10

# Sheet "fcr"
$ws = $wb.Worksheets.Item("fcr")
$ws.Cells.Item(2, 2).Value = 4.756275298391282
$ws.Cells.Item(2, 3).Value = 6.806138996138996
$ws.Cells.Item(2, 4).Value = 7.931085714285715
$ws.Cells.Item(3, 2).Value = 6.217189973614776
$ws.Cells.Item(3, 3).Value = 8.573076923076924
$ws.Cells.Item(3, 4).Value = 10.2555
$ws.Cells.Item(4, 2).Value = 7.284671232876713
$ws.Cells.Item(4, 3).Value = 12.22378571428571
$ws.Cells.Item(4, 4).Value = 8.577272727272728
$ws.Cells.Item(5, 2).Value = 4.427654830718414
$ws.Cells.Item(5, 3).Value = 6.140645161290323
$ws.Cells.Item(5, 4).Value = 7.556111111111111
$ws.Cells.Item(6, 2).Value = 9.176384976525821
$ws.Cells.Item(6, 3).Value = 11.96726315789474
$ws.Cells.Item(6, 4).Value = 23.533
$ws.Cells.Item(7, 2).Value = 5.545022421524664
$ws.Cells.Item(7, 3).Value = 7.683225806451613
$ws.Cells.Item(7, 4).Value = 3.79

# Sheet "rcr"
$ws = $wb.Worksheets.Item("rcr")
$ws.Cells.Item(2, 2).Value = 1.631124942369756
$ws.Cells.Item(2, 3).Value = 2.076743455497382
$ws.Cells.Item(2, 4).Value = 2.387987804878049
$ws.Cells.Item(3, 2).Value = 2.011741741741742
$ws.Cells.Item(3, 3).Value = 2.4452
$ws.Cells.Item(3, 4).Value = 2.6535
$ws.Cells.Item(4, 2).Value = 2.091850079744817
$ws.Cells.Item(4, 3).Value = 3.197404580152672
$ws.Cells.Item(4, 4).Value = 2.246451612903226
$ws.Cells.Item(5, 2).Value = 1.760020120724346
$ws.Cells.Item(5, 3).Value = 2.283059701492538
$ws.Cells.Item(5, 4).Value = 2.9146875
$ws.Cells.Item(6, 2).Value = 2.560963541666667
$ws.Cells.Item(6, 3).Value = 3.262613636363636
$ws.Cells.Item(6, 4).Value = 5.193
$ws.Cells.Item(7, 2).Value = 2.380109289617486
$ws.Cells.Item(7, 3).Value = 2.553793103448276
$ws.Cells.Item(7, 4).Value = 1.63
